{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two placeholder body paragraphs and the heading placeholder.\nconst firstParaIndex = items.findIndex(p => p.text === \"Put First Paragraph Here.\");\nconst secondParaIndex = items.findIndex(p => p.text === \"Put Secon Paragraph Here.\");\nconst headingIndex = items.findIndex(p => p.text === \"Put Title of Article Here\");\n\nif (firstParaIndex === -1 || secondParaIndex === -1 || headingIndex === -1) {\n  throw new Error(\"Could not locate expected placeholder paragraphs.\");\n}\n\nconst firstPara = items[firstParaIndex];\nconst secondPara = items[secondParaIndex];\nconst headingPara = items[headingIndex];\n\n// Replace the placeholder paragraph texts with the final copy.\nfirstPara.insertText(\n  \"This week, we will be taking a look at how we can fill the shape that we just created last week, by joining the two sides of our mirrored curve together. Then, after filling that shape, we will be creating a small hole in it.\",\n  \"Replace\"\n);\n\nsecondPara.insertText(\n  \"So, if you are interested in continuing on in this journey into turning those curves into shapes. Then please join us for our brand-new article this week entitled:\",\n  \"Replace\"\n);\n\n// Add the new single-space paragraph right after the second paragraph.\nsecondPara.insertParagraph(\" \", \"After\");\n\n// Update the article title heading.\nheadingPara.insertText(\"8 Fill and Cut Hole in Form\", \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the first placeholder paragraph's text.\n$r1 = $d.Content\n$r1.Find.Execute(\"Put First Paragraph Here.\", $false, $false, $false, $false, $false, $true, 1, $false, \"This week, we will be taking a look at how we can fill the shape that we just created last week, by joining the two sides of our mirrored curve together. Then, after filling that shape, we will be creating a small hole in it.\", 2)\n\n# Replace the second placeholder paragraph's text.\n$r2 = $d.Content\n$r2.Find.Execute(\"Put Secon Paragraph Here.\", $false, $false, $false, $false, $false, $true, 1, $false, \"So, if you are interested in continuing on in this journey into turning those curves into shapes. Then please join us for our brand-new article this week entitled:\", 2)\n\n# Insert a brand-new paragraph (containing a single space) right after the\n# second paragraph (paragraph index 3: Title=1, First=2, Second=3).\n$secondPara = $d.Paragraphs(3)\n$secondPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs(4)\n$newPara.Range.Text = \" \"\n\n# Replace the article-title heading placeholder's text.\n$r3 = $d.Content\n$r3.Find.Execute(\"Put Title of Article Here\", $false, $false, $false, $false, $false, $true, 1, $false, \"8 Fill and Cut Hole in Form\", 2)\n"}
